$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume values to match the latest scrape.
# D-column (Price) values are forced to Text format first so that Excel
# does not reinterpret numeric-looking strings (e.g. "569.75") as numbers,
# then the cell style is reset to Normal so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.484.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.181.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.21%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.179.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.23%  "
$ws.Range("E10").Value = "  -3.88%  "
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.732.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.32%  "
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.527.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.168.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "420.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("E21").Value = "  -3.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.489"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.18%  "
$ws.Range("E28").Value = "  -6.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -5.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("E35").Value = "  -4.01%  "
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  -5.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.728.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("E40").Value = "  -5.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "24.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("E44").Value = "  -7.12%  "
$ws.Range("E45").Value = "  -6.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.55%  "
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "294.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("E51").Value = "  -13.35%  "
